$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1029
$ws.Range("I2").Value = 1029
$ws.Range("K2").Value = 1029
$ws.Range("M2").Value = -916
$ws.Range("H17").Value = 926.15
$ws.Range("I17").Value = 300
$ws.Range("K17").Value = 900
$ws.Range("M17").Value = -732
$ws.Range("H40").Value = 10754.392
$ws.Range("I40").Value = 4778.5557
$ws.Range("J40").Value = 14596
$ws.Range("K40").Value = 4778.5557
$ws.Range("L40").Value = 14596
$ws.Range("M40").Value = -4603.5557
$ws.Range("N40").Value = -14946
$ws.Range("H41").Value = 1371
$ws.Range("I41").Value = 1725.625
$ws.Range("J41").Value = 1016.375
$ws.Range("K41").Value = 1725.625
$ws.Range("L41").Value = 1016.375
$ws.Range("M41").Value = -1285.625
$ws.Range("N41").Value = -1896.375
$ws.Range("H42").Value = 26.8
$ws.Range("J42").Value = 34
$ws.Range("L42").Value = 102
$ws.Range("N42").Value = -562
$ws.Range("H43").Value = 2798.6667
$ws.Range("I43").Value = 1899.5
$ws.Range("J43").Value = 3248.25
$ws.Range("K43").Value = 1899.5
$ws.Range("L43").Value = 3248.25
$ws.Range("M43").Value = -1830.5
$ws.Range("N43").Value = -3386.25
$ws.Range("H51").Value = 2833.3333
$ws.Range("I51").Value = 2250
$ws.Range("K51").Value = 2250
$ws.Range("M51").Value = -1766
$ws.Range("H74").Value = 5061.125
$ws.Range("I74").Value = 4824.909
$ws.Range("J74").Value = 5580.8
$ws.Range("K74").Value = 4824.909
$ws.Range("L74").Value = 5580.8
$ws.Range("M74").Value = -3888.909
$ws.Range("N74").Value = -7452.8
$ws.Range("H76").Value = 3152.95
$ws.Range("I76").Value = 2905.5
$ws.Range("J76").Value = 3524.125
$ws.Range("K76").Value = 2905.5
$ws.Range("L76").Value = 3524.125
$ws.Range("M76").Value = -2590.5
$ws.Range("N76").Value = -4154.125
$ws.Range("H77").Value = 5061.125
$ws.Range("I77").Value = 4824.909
$ws.Range("J77").Value = 5580.8
$ws.Range("K77").Value = 24124.545
$ws.Range("L77").Value = 27904
$ws.Range("M77").Value = -19444.545
$ws.Range("N77").Value = -37264
$ws.Range("H79").Value = 3152.95
$ws.Range("I79").Value = 2905.5
$ws.Range("J79").Value = 3524.125
$ws.Range("K79").Value = 2905.5
$ws.Range("L79").Value = 3524.125
$ws.Range("M79").Value = -1813.5
$ws.Range("N79").Value = -5708.125
$ws.Range("H82").Value = 8233
$ws.Range("I82").Value = 8233
$ws.Range("K82").Value = 24699
$ws.Range("M82").Value = -24293
$ws.Range("H85").Value = 8233
$ws.Range("I85").Value = 8233
$ws.Range("K85").Value = 24699
$ws.Range("M85").Value = -23295
$ws.Range("H86").Value = 3582.6191
$ws.Range("I86").Value = 3038
$ws.Range("J86").Value = 4308.778
$ws.Range("K86").Value = 3038
$ws.Range("L86").Value = 4308.778
$ws.Range("M86").Value = -1915
$ws.Range("N86").Value = -6554.778
$ws.Range("H89").Value = 3582.6191
$ws.Range("I89").Value = 3038
$ws.Range("J89").Value = 4308.778
$ws.Range("K89").Value = 15190
$ws.Range("L89").Value = 21543.89
$ws.Range("M89").Value = -9574
$ws.Range("N89").Value = -32775.89
$ws.Range("H106").Value = 4123
$ws.Range("I106").Value = 3938.5264
$ws.Range("J106").Value = 4999.25
$ws.Range("K106").Value = 3938.5264
$ws.Range("L106").Value = 4999.25
$ws.Range("M106").Value = -3307.5264
$ws.Range("N106").Value = -6261.25
$ws.Range("H111").Value = 669
$ws.Range("I111").Value = 406.66666
$ws.Range("J111").Value = 1062.5
$ws.Range("K111").Value = 1219.99998
$ws.Range("L111").Value = 3187.5
$ws.Range("M111").Value = 1847.00002
$ws.Range("N111").Value = -9321.5
$ws.Range("H134").Value = 40488.8
$ws.Range("J134").Value = 40488.8
$ws.Range("L134").Value = 40488.8
$ws.Range("N134").Value = -50628.8
$ws.Range("H138").Value = 1826.7273
$ws.Range("I138").Value = 1363.9231
$ws.Range("J138").Value = 2495.2222
$ws.Range("K138").Value = 4091.7693
$ws.Range("L138").Value = 7485.6666
$ws.Range("M138").Value = 1048.2307
$ws.Range("N138").Value = -17765.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 129998.336
$ws.Range("J23").Value = 129998.336
$ws.Range("L23").Value = 129998.336
$ws.Range("N23").Value = -130516.336
$ws.Range("H32").Value = 12868.621
$ws.Range("I32").Value = 5162
$ws.Range("J32").Value = 23786.334
$ws.Range("K32").Value = 5162
$ws.Range("L32").Value = 23786.334
$ws.Range("M32").Value = -4875
$ws.Range("N32").Value = -24360.334
$ws.Range("H74").Value = 35390.434
$ws.Range("I74").Value = 43100.543
$ws.Range("K74").Value = 43100.543
$ws.Range("M74").Value = -42226.543
$ws.Range("H77").Value = 35390.434
$ws.Range("I77").Value = 43100.543
$ws.Range("K77").Value = 215502.715
$ws.Range("M77").Value = -211134.715
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10695.286
$ws.Range("J20").Value = 1294
$ws.Range("L20").Value = 1294
$ws.Range("N20").Value = -1788
$ws.Range("H105").Value = 31423.588
$ws.Range("I105").Value = 35945.07
$ws.Range("J105").Value = 5199
$ws.Range("K105").Value = 35945.07
$ws.Range("L105").Value = 5199
$ws.Range("M105").Value = -34198.07
$ws.Range("N105").Value = -8693
$ws.Range("H107").Value = 14289659
$ws.Range("I107").Value = 16670518
$ws.Range("J107").Value = 4500
$ws.Range("K107").Value = 16670518
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = -16668598
$ws.Range("N107").Value = -8340
$ws.Range("H132").Value = 31181.25
$ws.Range("J132").Value = 31181.25
$ws.Range("L132").Value = 31181.25
$ws.Range("N132").Value = -41301.25
$ws.Range("H134").Value = 4102.5186
$ws.Range("I134").Value = 2626.3333
$ws.Range("K134").Value = 7878.999899999999
$ws.Range("M134").Value = -5343.999899999999
$ws.Range("H135").Value = 72741.14
$ws.Range("J135").Value = 72741.14
$ws.Range("L135").Value = 72741.14
$ws.Range("N135").Value = -82881.14
$ws.Range("H137").Value = 79040
$ws.Range("J137").Value = 79040
$ws.Range("L137").Value = 79040
$ws.Range("N137").Value = -89240
$ws.Range("H138").Value = 71076.086
$ws.Range("J138").Value = 71076.086
$ws.Range("L138").Value = 71076.086
$ws.Range("N138").Value = -81356.086

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 12000
$ws.Range("J32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -12632
$ws.Range("H107").Value = 1558.6666
$ws.Range("I107").Value = 1480.5
$ws.Range("K107").Value = 1480.5
$ws.Range("M107").Value = 439.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 499
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("H42").Value = 2809
$ws.Range("I42").Value = 1013.6667
$ws.Range("J42").Value = 5502
$ws.Range("K42").Value = 3041.0001
$ws.Range("L42").Value = 16506
$ws.Range("M42").Value = -2507.0001
$ws.Range("N42").Value = -17574
$ws.Range("H56").Value = 6518.75
$ws.Range("I56").Value = 6518.75
$ws.Range("K56").Value = 6518.75
$ws.Range("M56").Value = -5988.75
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 55091.418
$ws.Range("J70").Value = 14585.286
$ws.Range("L70").Value = 14585.286
$ws.Range("N70").Value = -15125.286
$ws.Range("H73").Value = 55091.418
$ws.Range("J73").Value = 14585.286
$ws.Range("L73").Value = 14585.286
$ws.Range("N73").Value = -16457.286
$ws.Range("H102").Value = 1604.7778
$ws.Range("I102").Value = 1604.7778
$ws.Range("K102").Value = 1604.7778
$ws.Range("M102").Value = 17.22219999999993
$ws.Range("H132").Value = 4153.8213
$ws.Range("I132").Value = 3467.238
$ws.Range("K132").Value = 10401.714
$ws.Range("M132").Value = -7871.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 10750
$ws.Range("I48").Value = 7000
$ws.Range("K48").Value = 7000
$ws.Range("M48").Value = -6339
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1504.56
$ws.Range("I107").Value = 1007.2941
$ws.Range("K107").Value = 3021.8823
$ws.Range("M107").Value = -1101.8823
$ws.Range("H123").Value = 74800
$ws.Range("J123").Value = 74800
$ws.Range("L123").Value = 74800
$ws.Range("N123").Value = -84600
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
